$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: bump the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update price list values in column D
$ws.Range("D33").Value = 698
$ws.Range("D34").Value = 899
$ws.Range("D35").Value = 970
$ws.Range("D36").Value = 1248
$ws.Range("D37").Value = 1578
$ws.Range("D38").Value = 1238
$ws.Range("D39").Value = 1565
$ws.Range("D40").Value = 1961
$ws.Range("D41").Value = 2277
